# Add season win/loss/tie record columns (Wins, Losses, Ties) to the
# roster sheet, per "Created functions to get season record".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look exactly like the existing header
# cells (bold, centered, bordered). Copy the formatting from the last
# existing header cell (AC1) onto the three new header cells, then set
# their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the team's 2003 season record: 83 wins, 79
# losses, 0 ties.
$lastRow = 55
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 83
    $ws.Cells.Item($row, 31).Value = 79
    $ws.Cells.Item($row, 32).Value = 0
}
